$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 43) holds the "Förändrad" date as a serial
# number. Bump each value by one day (45762 -> 45763).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45762) {
        $cell.Value = 45763
    }
}
